$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.829.24"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.92%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.093.25"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.98%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.705"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.95%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "245.54"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.16%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "53.80"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.61%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "59.11"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.65%  "
$ws.Range("E10").Value = "  -4.00%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0766"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.38%  "
$ws.Range("E12").Value = "  +1.13%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.914"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.51%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.85"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -8.64%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.400.45"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.10%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.47"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.95%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.140.12"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.22%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "36.806.90"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.13"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -8.33%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.96"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.26%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0879"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.44"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.34%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "239.64"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.97%  "
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("E25").Value = "  -3.68%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.63"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.42%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.15"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.45%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "166.63"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.85"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.128"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.88%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.21"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.16"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.60%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.71"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.43%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0607"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.58%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.44"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.44%  "
$ws.Range("E36").Value = "  +0.22%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.84"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.63%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0823"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.41%  "
$ws.Range("E39").Value = "  -5.47%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.16"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.33%  "
$ws.Range("B41").Value = "THORChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.90"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -8.03%  "
$ws.Range("B42").Value = "Cronos"
$ws.Range("C42").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0966"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.19%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0220"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.73%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "96.29"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.86"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.48%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.96"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +16.26%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.406.42"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +10.48%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "16.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -9.41%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.42"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.14%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.89"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.33%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.288.78"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.32%  "
